$wb = $excel.ActiveWorkbook

# ---- Step 1: insert the new 2022-Q3 sheet right after the summary sheet ----
$summarySheet = $wb.Worksheets.Item(1)
$templateSheet = $wb.Worksheets.Item("2022-Q2")
$templateSheet.Copy($null, $summarySheet)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# ---- Step 2: overwrite the copied data with the 2022-Q3 figures ----
$q3.Range("B2:G27").NumberFormat = "@"
$data = New-Object 'object[,]' 26,6
$data[0,0] = "162607"
$data[0,1] = "景顺长城资源垄断混合（LOF）"
$data[0,2] = "23.41"
$data[0,3] = "93.85"
$data[0,4] = "7.44"
$data[0,5] = "1.7417"
$data[1,0] = "001955"
$data[1,1] = "中欧养老产业混合A"
$data[1,2] = "22.51"
$data[1,3] = "92.70"
$data[1,4] = "7.07"
$data[1,5] = "1.5915"
$data[2,0] = "010429"
$data[2,1] = "中欧睿见混合A"
$data[2,2] = "18.74"
$data[2,3] = "91.22"
$data[2,4] = "7.17"
$data[2,5] = "1.3437"
$data[3,0] = "000772"
$data[3,1] = "景顺长城中国回报灵活配置混合"
$data[3,2] = "19.25"
$data[3,3] = "93.90"
$data[3,4] = "6.97"
$data[3,5] = "1.3417"
$data[4,0] = "166027"
$data[4,1] = "中欧创业板两年定期开放混合A"
$data[4,2] = "8.11"
$data[4,3] = "99.40"
$data[4,4] = "10.07"
$data[4,5] = "0.8167"
$data[5,0] = "011710"
$data[5,1] = "中欧睿泽混合A"
$data[5,2] = "8.33"
$data[5,3] = "90.88"
$data[5,4] = "7.95"
$data[5,5] = "0.6622"
$data[6,0] = "001487"
$data[6,1] = "宝盈优势产业灵活配置混合A"
$data[6,2] = "10.11"
$data[6,3] = "91.85"
$data[6,4] = "2.85"
$data[6,5] = "0.2881"
$data[7,0] = "007202"
$data[7,1] = "天弘优质成长企业精选混合A"
$data[7,2] = "5.80"
$data[7,3] = "93.00"
$data[7,4] = "4.89"
$data[7,5] = "0.2836"
$data[8,0] = "009791"
$data[8,1] = "中欧创业板两年定期开放混合C"
$data[8,2] = "2.81"
$data[8,3] = "99.40"
$data[8,4] = "10.07"
$data[8,5] = "0.2830"
$data[9,0] = "012778"
$data[9,1] = "中欧养老产业混合C"
$data[9,2] = "2.80"
$data[9,3] = "92.70"
$data[9,4] = "7.07"
$data[9,5] = "0.1980"
$data[10,0] = "420005"
$data[10,1] = "天弘周期策略混合A"
$data[10,2] = "3.10"
$data[10,3] = "93.08"
$data[10,4] = "6.18"
$data[10,5] = "0.1916"
$data[11,0] = "001075"
$data[11,1] = "宝盈转型动力灵活配置混合A"
$data[11,2] = "4.35"
$data[11,3] = "91.90"
$data[11,4] = "4.17"
$data[11,5] = "0.1814"
$data[12,0] = "000586"
$data[12,1] = "景顺长城中小创精选股票"
$data[12,2] = "2.21"
$data[12,3] = "93.50"
$data[12,4] = "6.87"
$data[12,5] = "0.1518"
$data[13,0] = "012771"
$data[13,1] = "宝盈优势产业灵活配置混合C"
$data[13,2] = "3.62"
$data[13,3] = "91.85"
$data[13,4] = "2.85"
$data[13,5] = "0.1032"
$data[14,0] = "011711"
$data[14,1] = "中欧睿泽混合C"
$data[14,2] = "0.86"
$data[14,3] = "90.88"
$data[14,4] = "7.95"
$data[14,5] = "0.0684"
$data[15,0] = "260115"
$data[15,1] = "景顺长城中小盘混合"
$data[15,2] = "0.92"
$data[15,3] = "92.87"
$data[15,4] = "5.18"
$data[15,5] = "0.0477"
$data[16,0] = "015458"
$data[16,1] = "天弘周期策略混合C"
$data[16,2] = "0.68"
$data[16,3] = "93.08"
$data[16,4] = "6.18"
$data[16,5] = "0.0420"
$data[17,0] = "015769"
$data[17,1] = "天弘低碳经济混合A"
$data[17,2] = "1.19"
$data[17,3] = "79.16"
$data[17,4] = "3.48"
$data[17,5] = "0.0414"
$data[18,0] = "015770"
$data[18,1] = "天弘低碳经济混合C"
$data[18,2] = "1.03"
$data[18,3] = "79.16"
$data[18,4] = "3.48"
$data[18,5] = "0.0358"
$data[19,0] = "000573"
$data[19,1] = "天弘通利混合"
$data[19,2] = "1.01"
$data[19,3] = "79.25"
$data[19,4] = "3.20"
$data[19,5] = "0.0323"
$data[20,0] = "004694"
$data[20,1] = "天弘策略精选灵活配置混合A"
$data[20,2] = "0.80"
$data[20,3] = "86.39"
$data[20,4] = "3.68"
$data[20,5] = "0.0294"
$data[21,0] = "015389"
$data[21,1] = "宝盈转型动力灵活配置混合C"
$data[21,2] = "0.24"
$data[21,3] = "91.90"
$data[21,4] = "4.17"
$data[21,5] = "0.0100"
$data[22,0] = "015481"
$data[22,1] = "中欧睿见混合C"
$data[22,2] = "0.08"
$data[22,3] = "91.22"
$data[22,4] = "7.17"
$data[22,5] = "0.0057"
$data[23,0] = "007084"
$data[23,1] = "天治转型升级混合"
$data[23,2] = "0.11"
$data[23,3] = "91.86"
$data[23,4] = "3.02"
$data[23,5] = "0.0033"
$data[24,0] = "015460"
$data[24,1] = "天弘优质成长企业精选混合C"
$data[24,2] = "0.06"
$data[24,3] = "93.00"
$data[24,4] = "4.89"
$data[24,5] = "0.0029"
$data[25,0] = "004748"
$data[25,1] = "天弘策略精选灵活配置混合C"
$data[25,2] = "0.06"
$data[25,3] = "86.39"
$data[25,4] = "3.68"
$data[25,5] = "0.0022"
$q3.Range("B2:G27").Value = $data

$hData = New-Object 'object[,]' 26,1
$hData[0,0] = 9
$hData[1,0] = 10
$hData[2,0] = 8
$hData[3,0] = 9
$hData[4,0] = 2
$hData[5,0] = 8
$hData[6,0] = 8
$hData[7,0] = 5
$hData[8,0] = 2
$hData[9,0] = 10
$hData[10,0] = 7
$hData[11,0] = 2
$hData[12,0] = 7
$hData[13,0] = 8
$hData[14,0] = 8
$hData[15,0] = 8
$hData[16,0] = 7
$hData[17,0] = 5
$hData[18,0] = 5
$hData[19,0] = 5
$hData[20,0] = 2
$hData[21,0] = 2
$hData[22,0] = 8
$hData[23,0] = 7
$hData[24,0] = 5
$hData[25,0] = 2
$q3.Range("H2:H27").Value = $hData

# ---- Step 3: update the summary (总计) sheet with the new quarter row ----
$summarySheet.Range("A2:D8").Copy($summarySheet.Range("A3:D9"))
$summarySheet.Range("B2").Value = "2022-Q3"
$summarySheet.Range("C2").Value = 26
$summarySheet.Range("D2").Value = 9.5
$aData = New-Object 'object[,]' 8,1
$aData[0,0] = 0
$aData[1,0] = 1
$aData[2,0] = 2
$aData[3,0] = 3
$aData[4,0] = 4
$aData[5,0] = 5
$aData[6,0] = 6
$aData[7,0] = 7
$summarySheet.Range("A2:A9").Value = $aData
